$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2024291497975708
$ws.Range("C2").Value = 0.5465587044534413
$ws.Range("J2").Value = 0.02024291497975709
$ws.Range("P2").Value = 0.1295546558704453
$ws.Range("S2").Value = 0.1012145748987854
$ws.Range("C3").Value = 0.03623188405797102
$ws.Range("J3").Value = 0.01449275362318841
$ws.Range("P3").Value = 0.6956521739130435
$ws.Range("S3").Value = 0.2536231884057971
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.7826086956521739
$ws.Range("S4").Value = 0.1956521739130435
$ws.Range("B6").Value = 0.04950495049504951
$ws.Range("D6").Value = 0.01485148514851485
$ws.Range("F6").Value = 0.05445544554455446
$ws.Range("J6").Value = 0.1683168316831683
$ws.Range("O6").Value = 0.0297029702970297
$ws.Range("Q6").Value = 0.1138613861386139
$ws.Range("R6").Value = 0.1089108910891089
$ws.Range("S6").Value = 0.4603960396039604
$ws.Range("B7").Value = 0.07291666666666667
$ws.Range("D7").Value = 0.03645833333333334
$ws.Range("F7").Value = 0.046875
$ws.Range("J7").Value = 0.109375
$ws.Range("Q7").Value = 0.234375
$ws.Range("R7").Value = 0.046875
$ws.Range("S7").Value = 0.453125
$ws.Range("B8").Value = 0.08620689655172414
$ws.Range("D8").Value = 0.01939655172413793
$ws.Range("E8").Value = 0.002155172413793103
$ws.Range("F8").Value = 0.08405172413793104
$ws.Range("J8").Value = 0.07974137931034483
$ws.Range("O8").Value = 0.03017241379310345
$ws.Range("Q8").Value = 0.1918103448275862
$ws.Range("R8").Value = 0.08189655172413793
$ws.Range("S8").Value = 0.4245689655172414
$ws.Range("B9").Value = 0.06018518518518518
$ws.Range("D9").Value = 0.01388888888888889
$ws.Range("F9").Value = 0.06481481481481481
$ws.Range("J9").Value = 0.1064814814814815
$ws.Range("O9").Value = 0.01388888888888889
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.4907407407407408
$ws.Range("B10").Value = 0.1091417910447761
$ws.Range("D10").Value = 0.02332089552238806
$ws.Range("F10").Value = 0.06996268656716417
$ws.Range("J10").Value = 0.105410447761194
$ws.Range("O10").Value = 0.03078358208955224
$ws.Range("Q10").Value = 0.1996268656716418
$ws.Range("R10").Value = 0.08582089552238806
$ws.Range("S10").Value = 0.3759328358208955
$ws.Range("G11").Value = 0.133587786259542
$ws.Range("J11").Value = 0.07251908396946564
$ws.Range("K11").Value = 0.1908396946564886
$ws.Range("L11").Value = 0.5916030534351145
$ws.Range("S11").Value = 0.01145038167938931
$ws.Range("G12").Value = 0.7935483870967742
$ws.Range("J12").Value = 0.167741935483871
$ws.Range("L12").Value = 0.01935483870967742
$ws.Range("S12").Value = 0.01935483870967742
$ws.Range("G13").Value = 0.7037037037037037
$ws.Range("J13").Value = 0.2962962962962963
$ws.Range("F15").Value = 0.01015228426395939
$ws.Range("H15").Value = 0.16751269035533
$ws.Range("I15").Value = 0.07106598984771574
$ws.Range("J15").Value = 0.2944162436548223
$ws.Range("K15").Value = 0.06598984771573604
$ws.Range("M15").Value = 0.01015228426395939
$ws.Range("O15").Value = 0.05583756345177665
$ws.Range("S15").Value = 0.3248730964467005
$ws.Range("F16").Value = 0.01875
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.08125
$ws.Range("J16").Value = 0.43125
$ws.Range("K16").Value = 0.075
$ws.Range("M16").Value = 0.01875
$ws.Range("O16").Value = 0.04375
$ws.Range("S16").Value = 0.13125
$ws.Range("F17").Value = 0.01231527093596059
$ws.Range("H17").Value = 0.2142857142857143
$ws.Range("I17").Value = 0.1009852216748768
$ws.Range("J17").Value = 0.416256157635468
$ws.Range("M17").Value = 0.0270935960591133
$ws.Range("N17").Value = 0.004926108374384237
$ws.Range("O17").Value = 0.03940886699507389
$ws.Range("S17").Value = 0.1133004926108374
$ws.Range("F18").Value = 0.005681818181818182
$ws.Range("H18").Value = 0.2159090909090909
$ws.Range("I18").Value = 0.09659090909090909
$ws.Range("J18").Value = 0.3920454545454545
$ws.Range("K18").Value = 0.1136363636363636
$ws.Range("M18").Value = 0.01136363636363636
$ws.Range("O18").Value = 0.06818181818181818
$ws.Range("S18").Value = 0.09659090909090909
$ws.Range("F19").Value = 0.01341156747694887
$ws.Range("H19").Value = 0.2321877619446773
$ws.Range("I19").Value = 0.1123218776194468
$ws.Range("J19").Value = 0.3528918692372171
$ws.Range("K19").Value = 0.1123218776194468
$ws.Range("M19").Value = 0.03017602682313495
$ws.Range("N19").Value = 0.001676445934618609
$ws.Range("O19").Value = 0.06035205364626991
$ws.Range("S19").Value = 0.08466051969823973
